$wb = $excel.ActiveWorkbook

# Add "Department1" sheet after the last existing sheet ("Regional Office Locations")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Department1"
$ws3.Range("A1").Value = "Correctional Services Oversight & Investigations"
$ws3.Range("A2").Value = "Family Liaison Support"
$ws3.Range("A3").Value = "Institution Admin Support"
$ws3.Range("A4").Value = "Operational Support Division"
$ws3.Range("A5").Value = "Safer Team"
$ws3.Range("A6").Value = "Statistical Analysis Unit"

# Add "Department3" sheet after "Department1"
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Department3"
$ws4.Range("A1").Value = "ADM Office of the Institutional Services"
$ws4.Range("A2").Value = "Executive Directors Office"
$ws4.Range("A3").Value = "Information Management Unit"

# Department3 ends up as the newly active / selected sheet (matches activeTab="3" in workbook.xml)
$ws4.Activate()
